$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 113; existing rows 113:126 shift down to 114:127.
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new weekly price record
# (Feria Lagunitas de Puerto Montt - Apio, Americana (o), Primera).
$ws.Cells.Item(113, 1).Value = 4
$ws.Cells.Item(113, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(113, 3).Value = "Los Lagos"
$ws.Cells.Item(113, 4).Value = 44449
$ws.Cells.Item(113, 5).Value = 10
$ws.Cells.Item(113, 6).Value = 100112017
$ws.Cells.Item(113, 7).Value = "Apio"
$ws.Cells.Item(113, 8).Value = "Americana (o)"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 40
$ws.Cells.Item(113, 11).Value = 12000
$ws.Cells.Item(113, 12).Value = 12000
$ws.Cells.Item(113, 13).Value = 12000
$ws.Cells.Item(113, 14).Value = "`$/docena de matas"
$ws.Cells.Item(113, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(113, 16).Value = 2000
$ws.Cells.Item(113, 17).Value = 6
$ws.Cells.Item(113, 18).Value = "Hortaliza"
